# Scheduled-runner refresh of market-price-derived columns (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the ALC/ARM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
# Values below were recomputed upstream; this just writes the refreshed figures
# into the existing Leve tables cell by cell (adding/clearing cells where a
# profit figure newly appears or a LeveProfitHQ stops applying).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 694.4167
$ws.Range("I33").Value = 272.83334
$ws.Range("K33").Value = 272.83334
$ws.Range("M33").Value = -43.83334000000002
$ws.Range("H92").Value = 2703.111
$ws.Range("I92").Value = 4108.25
$ws.Range("J92").Value = 1579
$ws.Range("K92").Value = 4108.25
$ws.Range("L92").Value = 1579
$ws.Range("M92").Value = -2860.25
$ws.Range("N92").Value = -4075
$ws.Range("H111").Value = 7912.4287
$ws.Range("I111").Value = 8000
$ws.Range("J111").Value = 7897.8335
$ws.Range("K111").Value = 24000
$ws.Range("L111").Value = 23693.5005
$ws.Range("M111").Value = -20933
$ws.Range("N111").Value = -29827.5005
$ws.Range("H112").Value = 1689.8286
$ws.Range("J112").Value = 1705.5938
$ws.Range("L112").Value = 5116.7814
$ws.Range("N112").Value = -7332.7814
$ws.Range("H132").Value = 23556.924
$ws.Range("I132").Value = 3987.4707
$ws.Range("K132").Value = 11962.4121
$ws.Range("M132").Value = -9432.4121
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1753.75
$ws.Range("I21").Value = 1171.6666
$ws.Range("J21").Value = 3500
$ws.Range("K21").Value = 1171.6666
$ws.Range("L21").Value = 3500
$ws.Range("M21").Value = -797.6666
$ws.Range("N21").Value = -4248
$ws.Range("H38").Value = 779
$ws.Range("I38").Value = 779
$ws.Range("K38").Value = 779
$ws.Range("M38").Value = -312
$ws.Range("H96").Value = 39498.75
$ws.Range("J96").Value = 39498.75
$ws.Range("L96").Value = 39498.75
$ws.Range("N96").Value = -44990.75
$ws.Range("H132").Value = 2650.158
$ws.Range("I132").Value = 2594.2942
$ws.Range("K132").Value = 7782.882599999999
$ws.Range("M132").Value = -5252.882599999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 234.5
$ws.Range("I2").Value = 400
$ws.Range("J2").Value = 179.33333
$ws.Range("K2").Value = 400
$ws.Range("L2").Value = 179.33333
$ws.Range("M2").Value = -287
$ws.Range("N2").Value = -405.33333
$ws.Range("H5").Value = 1195.5
$ws.Range("I5").Value = 1500
$ws.Range("J5").Value = 1043.25
$ws.Range("K5").Value = 1500
$ws.Range("L5").Value = 1043.25
$ws.Range("M5").Value = -1388
$ws.Range("N5").Value = -1267.25
$ws.Range("H7").Value = 168.8125
$ws.Range("I7").Value = 24.166666
$ws.Range("J7").Value = 255.6
$ws.Range("K7").Value = 24.166666
$ws.Range("L7").Value = 255.6
$ws.Range("M7").Value = 88.83333400000001
$ws.Range("N7").Value = -481.6
$ws.Range("H11").Value = 7166.3335
$ws.Range("J11").Value = 7166.3335
$ws.Range("L11").Value = 7166.3335
$ws.Range("N11").Value = -7446.3335
$ws.Range("H12").Value = 3899.2
$ws.Range("J12").Value = 5332
$ws.Range("L12").Value = 5332
$ws.Range("N12").Value = -5672
$ws.Range("H38").Value = 15897.25
$ws.Range("I38").Value = 15897.25
$ws.Range("K38").Value = 15897.25
$ws.Range("M38").Value = -15520.25
$ws.Range("H45").Value = 8000
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H46").Value = 15897.25
$ws.Range("I46").Value = 15897.25
$ws.Range("K46").Value = 15897.25
$ws.Range("M46").Value = -15686.25
$ws.Range("H107").Value = 5953306
$ws.Range("I107").Value = 8929440
$ws.Range("J107").Value = 1037.625
$ws.Range("K107").Value = 8929440
$ws.Range("L107").Value = 1037.625
$ws.Range("M107").Value = -8927520
$ws.Range("N107").Value = -4877.625
$ws.Range("H109").Value = 19175.54
$ws.Range("J109").Value = 19175.54
$ws.Range("L109").Value = 19175.54
$ws.Range("N109").Value = -21255.54
$ws.Range("H132").Value = 3330.7058
$ws.Range("I132").Value = 2693.3635
$ws.Range("K132").Value = 8080.0905
$ws.Range("M132").Value = -5550.0905
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1161.8889
$ws.Range("I7").Value = 1454.6364
$ws.Range("J7").Value = 701.8570999999999
$ws.Range("K7").Value = 4363.9092
$ws.Range("L7").Value = 2105.5713
$ws.Range("M7").Value = -4251.9092
$ws.Range("N7").Value = -2329.5713
$ws.Range("H29").Value = 759.8
$ws.Range("I29").Value = 885.4286
$ws.Range("J29").Value = 466.66666
$ws.Range("K29").Value = 2656.2858
$ws.Range("L29").Value = 1399.99998
$ws.Range("M29").Value = -2379.2858
$ws.Range("N29").Value = -1953.99998
$ws.Range("H50").Value = 788.6667
$ws.Range("J50").Value = 1121.6666
$ws.Range("L50").Value = 3364.9998
$ws.Range("N50").Value = -4326.9998
$ws.Range("H53").Value = 788.6667
$ws.Range("J53").Value = 1121.6666
$ws.Range("L53").Value = 3364.9998
$ws.Range("N53").Value = -4326.9998
$ws.Range("H82").Value = 16667499
$ws.Range("I82").Value = 16667499
$ws.Range("K82").Value = 50002497
$ws.Range("M82").Value = -50002091
$ws.Range("H85").Value = 16667499
$ws.Range("I85").Value = 16667499
$ws.Range("K85").Value = 50002497
$ws.Range("M85").Value = -50001093
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 30441.8
$ws.Range("J94").Value = 30441.8
$ws.Range("L94").Value = 30441.8
$ws.Range("N94").Value = -31793.8
$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -56884
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 1355.3334
$ws.Range("I32").Value = 1518.5
$ws.Range("J32").Value = 50
$ws.Range("K32").Value = 1518.5
$ws.Range("L32").Value = 50
$ws.Range("M32").Value = -1201.5
$ws.Range("N32").Value = -684
$ws.Range("H40").Value = 2542.2
$ws.Range("I40").Value = 2107.6365
$ws.Range("J40").Value = 3737.25
$ws.Range("K40").Value = 2107.6365
$ws.Range("L40").Value = 3737.25
$ws.Range("M40").Value = -1971.6365
$ws.Range("N40").Value = -4009.25
$ws.Range("H55").Value = 862.5
$ws.Range("I55").Value = 996
$ws.Range("J55").Value = 462
$ws.Range("K55").Value = 996
$ws.Range("L55").Value = 462
$ws.Range("M55").Value = -823
$ws.Range("N55").Value = -808
$ws.Range("H132").Value = 2368.524
$ws.Range("I132").Value = 1872.5883
$ws.Range("K132").Value = 5617.7649
$ws.Range("M132").Value = -3087.7649
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 8271.4
$ws.Range("I74").Value = 6400
$ws.Range("J74").Value = 8739.25
$ws.Range("K74").Value = 6400
$ws.Range("L74").Value = 8739.25
$ws.Range("M74").Value = -5464
$ws.Range("N74").Value = -10611.25
$ws.Range("H77").Value = 8271.4
$ws.Range("I77").Value = 6400
$ws.Range("J77").Value = 8739.25
$ws.Range("K77").Value = 19200
$ws.Range("L77").Value = 26217.75
$ws.Range("M77").Value = -14520
$ws.Range("N77").Value = -35577.75
$ws.Range("H109").Value = 20650
$ws.Range("J109").Value = 20650
$ws.Range("L109").Value = 20650
$ws.Range("N109").Value = -23424
$ws.Range("H132").Value = 1328.9546
$ws.Range("J132").Value = 3410
$ws.Range("L132").Value = 10230
$ws.Range("N132").Value = -15290
$ws.Range("H141").Value = 82222.22
$ws.Range("J141").Value = 82222.22
$ws.Range("L141").Value = 82222.22
$ws.Range("N141").Value = -92582.22
